$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared across all rows 144-163 in this subset (Vega Monumental Concepcion - Perejil)
$constA = 11
$constB = "Vega Monumental Concepción"
$constC = "Bíobío"
$constE = 8
$constF = 100112044
$constG = "Perejil"
$constH = "Sin especificar"
$constN = "$/atado 0,5 a 1 kilo"
$constO = "Región de Ñuble"
$constQ = 1
$constR = "Hortaliza"

# Row data: row number, Fecha (D), Calidad (I), Volumen (J), Precio minimo (K), Precio maximo (L), Precio promedio (M), Precio $/Kg (P)
$rows = @(
    ,@(144, 44783, "Primera", 200, 700, 800, 750, 750)
    ,@(145, 44783, "Segunda", 100, 600, 600, 600, 600)
    ,@(146, 44425, "Primera", 200, 600, 700, 650, 650)
    ,@(147, 44425, "Segunda", 100, 500, 500, 500, 500)
    ,@(148, 44512, "Primera", 200, 600, 700, 650, 650)
    ,@(149, 44512, "Segunda", 100, 500, 500, 500, 500)
    ,@(150, 44285, "Primera", 200, 600, 700, 650, 650)
    ,@(151, 44285, "Segunda", 100, 500, 500, 500, 500)
    ,@(152, 44362, "Primera", 200, 600, 700, 650, 650)
    ,@(153, 44362, "Segunda", 100, 500, 500, 500, 500)
    ,@(154, 44355, "Primera", 200, 600, 700, 650, 650)
    ,@(155, 44355, "Segunda", 100, 500, 500, 500, 500)
    ,@(156, 44391, "Primera", 200, 600, 700, 650, 650)
    ,@(157, 44391, "Segunda", 100, 500, 500, 500, 500)
    ,@(158, 44453, "Primera", 200, 600, 700, 650, 650)
    ,@(159, 44453, "Segunda", 100, 500, 500, 500, 500)
    ,@(160, 44609, "Primera", 200, 600, 700, 650, 650)
    ,@(161, 44609, "Segunda", 100, 500, 500, 500, 500)
    ,@(162, 44358, "Primera", 200, 600, 700, 650, 650)
    ,@(163, 44358, "Segunda", 100, 500, 500, 500, 500)
)

foreach ($row in $rows) {
    $r = $row[0]
    $d = $row[1]
    $qual = $row[2]
    $vol = $row[3]
    $kmin = $row[4]
    $lmax = $row[5]
    $mavg = $row[6]
    $pkg = $row[7]

    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $qual
    $ws.Cells.Item($r, 10).Value = $vol
    $ws.Cells.Item($r, 11).Value = $kmin
    $ws.Cells.Item($r, 12).Value = $lmax
    $ws.Cells.Item($r, 13).Value = $mavg
    $ws.Cells.Item($r, 14).Value = $constN
    $ws.Cells.Item($r, 15).Value = $constO
    $ws.Cells.Item($r, 16).Value = $pkg
    $ws.Cells.Item($r, 17).Value = $constQ
    $ws.Cells.Item($r, 18).Value = $constR
}

Write-Output "Done updating rows 144-163"